$d = $word.ActiveDocument

# Locate the paragraphs to remove:
#   - the blank paragraph right before "Ver no Jupiter Salvar em pdf Salvar em docx"
#   - the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph itself
#   - the "© 2020 . Contact: luizeleno@usp.br. ..." paragraph
# by scanning paragraph text, so the script does not depend on brittle hardcoded
# paragraph indices.
$count = $d.Paragraphs.Count
$startIdx = -1
$endIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Ver no Jupiter*") {
        $startIdx = $i - 1
    }
    if ($t -like "*Contact: luizeleno*") {
        $endIdx = $i
    }
}

if ($startIdx -gt 0 -and $endIdx -ge $startIdx) {
    $start = $d.Paragraphs.Item($startIdx).Range.Start
    $end = $d.Paragraphs.Item($endIdx).Range.End
    $r = $d.Range($start, $end)
    $r.Delete()
}
